$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the sheet-tab area in the window (cosmetic window setting).
$excel.ActiveWindow.TabRatio = 0.991

# Fill in the new header cells (D1:H1), preserving the existing styles on D1/E1/F1
# and using the default style for the newly introduced G1/H1 columns.
$ws.Range("D1").Value = "End Date"
$ws.Range("E1").Value = "Business"
$ws.Range("F1").Value = "Status"
$ws.Range("G1").Value = "Database"
$ws.Range("H1").Value = "Report"

# Normalize row heights for the data rows (13.8 -> 15).
for ($r = 1; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# Move the active selection as recorded after the edit.
$ws.Range("E11").Select()
